$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the CSE 208 / section 1 exam schedule cell (E6):
# venue stays "Room PB - 105 (Permanent Campus)" but the date/time moves to
# 6.1.2020 (Monday) 10:45 - 12:45.
$header = "Final Exam Schedule"
$body = "Venue: Room PB - 105 (Permanent Campus)" + "`n" + "Date: 6.1.2020 (Monday)" + "`n" + "Time: 10:45 - 12:45"
$full = $header + "`n" + $body

$cell = $ws.Range("E6")
$cell.Value = $full

# Preserve the original rich-text formatting: bold "Final Exam Schedule"
# heading followed by a regular-weight details block.
$headerLen = $header.Length
$cell.Characters(1, $headerLen).Font.Bold = $true
$cell.Characters(1, $headerLen).Font.Size = 11
$cell.Characters(1, $headerLen).Font.Name = "Calibri"

$restStart = $headerLen + 1
$restLen = $full.Length - $headerLen
$cell.Characters($restStart, $restLen).Font.Bold = $false
$cell.Characters($restStart, $restLen).Font.Size = 11
$cell.Characters($restStart, $restLen).Font.Name = "Calibri"

# Match the saved workbook's updated cursor position.
$ws.Range("D6").Select()
